$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 781 (shifts the existing rows 781..822 down to 782..823)
$ws.Rows(781).Insert()

# Fill in the new row's data.
# Column A holds a date written as plain text (e.g. "2026/12/29"), not a
# real Excel date value, so use a quote-prefixed Formula assignment (which
# strips the leading apostrophe but keeps the value as text) and then reset
# the style back to Normal so no stray number-format is left on the cell.
$ws.Range("A781").Formula = "'2026/02/11"
$ws.Range("A781").Style = "Normal"
$ws.Range("B781").Value = "水"
$ws.Range("C781").Value = 17
$ws.Range("D781").Value = 167
